# "fix the merge error" - Player.xlsx / Property sheet
#
# A previous merge had left several "View" (column F) flags un-set and a
# couple of Private/Save/View flags mixed up on rows 76-77, plus row 78's
# Public flag incorrectly left on. This restores the intended values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property")
$ws.Activate()

# --- header row: F1 re-affirmed as "View" (column header text) ---
$ws.Range("F1").Value = "View"

# --- rows 68-75: these rows were missing the "View" (column F) flag
#     entirely; the merge should have carried it over as TRUE like the
#     surrounding rows. ---
$ws.Range("F68").Value = $true
$ws.Range("F69").Value = $true
$ws.Range("F70").Value = $true
$ws.Range("F71").Value = $true
$ws.Range("F72").Value = $true
$ws.Range("F73").Value = $true
$ws.Range("F74").Value = $true
$ws.Range("F75").Value = $true

# --- row 76 (GameID): Private/Save should be FALSE, View should be TRUE ---
$ws.Range("D76").Value = $false
$ws.Range("E76").Value = $false
$ws.Range("F76").Value = $true

# --- row 77 (GateID): Private/Save should be FALSE, View should be TRUE ---
$ws.Range("D77").Value = $false
$ws.Range("E77").Value = $false
$ws.Range("F77").Value = $true

# --- row 78 (GuildID): Public flag should be FALSE ---
$ws.Range("C78").Value = $false

# --- restore the scroll position / selection recorded for the sheet ---
$win = $excel.ActiveWindow
$win.ScrollRow = 37
$win.ScrollColumn = 1
$ws.Range("C78").Select()
